# Refresh Universalis market-price snapshots (current/NQ/HQ prices + derived leve
# profit columns H:N) across all eight crafter job sheets, per the scheduled-runner
# data pull. Only data cells change; layout/tables/styles are untouched.
# Note: a couple of rows gain/lose an M or N cell entirely (profit column only
# applies when that price tier - NQ vs HQ - is actually in play for the leve).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1042.3
$ws.Range("I12").Value = 1302.3334
$ws.Range("J12").Value = 652.25
$ws.Range("K12").Value = 1302.3334
$ws.Range("L12").Value = 652.25
$ws.Range("M12").Value = -1132.3334
$ws.Range("N12").Value = -992.25
$ws.Range("H17").Value = 650
$ws.Range("I17").Value = 650
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1950
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1782
$ws.Range("N17").ClearContents()
$ws.Range("H39").Value = 724.1818
$ws.Range("I39").Value = 347
$ws.Range("J39").Value = 1730
$ws.Range("K39").Value = 1041
$ws.Range("L39").Value = 5190
$ws.Range("M39").Value = -745
$ws.Range("N39").Value = -5782
$ws.Range("H40").Value = 3499.4
$ws.Range("J40").Value = 3574.25
$ws.Range("L40").Value = 3574.25
$ws.Range("N40").Value = -3924.25
$ws.Range("H62").Value = 3944
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 3944
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H137").Value = 1260.0741
$ws.Range("J137").Value = 1203.2858
$ws.Range("L137").Value = 3609.8574
$ws.Range("N137").Value = -8709.857400000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 815.75
$ws.Range("I61").Value = 503.7143
$ws.Range("K61").Value = 503.7143
$ws.Range("M61").Value = -291.7143
$ws.Range("H74").Value = 793.4286
$ws.Range("I74").Value = 793.4286
$ws.Range("K74").Value = 793.4286
$ws.Range("M74").Value = 80.57140000000004
$ws.Range("H77").Value = 793.4286
$ws.Range("I77").Value = 793.4286
$ws.Range("K77").Value = 3967.143
$ws.Range("M77").Value = 400.857
$ws.Range("H102").Value = 3069.6
$ws.Range("I102").Value = 3299.5557
$ws.Range("K102").Value = 3299.5557
$ws.Range("M102").Value = -1677.5557
$ws.Range("H136").Value = 815.75
$ws.Range("I136").Value = 503.7143
$ws.Range("K136").Value = 1511.1429
$ws.Range("M136").Value = 1038.8571

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 7987.8335
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -10492
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2389.8
$ws.Range("I31").Value = 1425
$ws.Range("J31").Value = 3033
$ws.Range("K31").Value = 1425
$ws.Range("L31").Value = 3033
$ws.Range("M31").Value = -1130
$ws.Range("N31").Value = -3623
$ws.Range("H34").Value = 2389.8
$ws.Range("I34").Value = 1425
$ws.Range("J34").Value = 3033
$ws.Range("K34").Value = 1425
$ws.Range("L34").Value = 3033
$ws.Range("M34").Value = -1223
$ws.Range("N34").Value = -3437
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232
$ws.Range("H99").Value = 6191.769
$ws.Range("J99").Value = 5498.5
$ws.Range("L99").Value = 5498.5
$ws.Range("N99").Value = -8494.5
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -20242
$ws.Range("H107").Value = 435.83334
$ws.Range("I107").Value = 390.8889
$ws.Range("K107").Value = 390.8889
$ws.Range("M107").Value = 1529.1111
$ws.Range("H122").Value = 2707.5
$ws.Range("I122").Value = 2721.4285
$ws.Range("K122").Value = 8164.2855
$ws.Range("M122").Value = -5714.2855
$ws.Range("H126").Value = 6191.769
$ws.Range("J126").Value = 5498.5
$ws.Range("L126").Value = 16495.5
$ws.Range("N126").Value = -21435.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 94.75
$ws.Range("I33").Value = 93
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 558
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -275
$ws.Range("N33").Value = -1166
$ws.Range("H113").Value = 1216.5555
$ws.Range("J113").Value = 1514.8334
$ws.Range("L113").Value = 4544.5002
$ws.Range("N113").Value = -8884.5002
$ws.Range("H131").Value = 984.5217
$ws.Range("J131").Value = 989.7619
$ws.Range("L131").Value = 2969.2857
$ws.Range("N131").Value = -13049.2857

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 752.8182
$ws.Range("I97").Value = 375.66666
$ws.Range("J97").Value = 2450
$ws.Range("K97").Value = 375.66666
$ws.Range("L97").Value = 2450
$ws.Range("M97").Value = 120.33334
$ws.Range("N97").Value = -3442
$ws.Range("H102").Value = 3331.111
$ws.Range("I102").Value = 3326
$ws.Range("K102").Value = 3326
$ws.Range("M102").Value = -1704
$ws.Range("H113").Value = 1198.5
$ws.Range("I113").Value = 1198.5
$ws.Range("K113").Value = 1198.5
$ws.Range("M113").Value = 971.5
$ws.Range("H122").Value = 1690.5
$ws.Range("I122").Value = 1463.8125
$ws.Range("K122").Value = 4391.4375
$ws.Range("M122").Value = -1941.4375

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 560.4375
$ws.Range("I55").Value = 167.33333
$ws.Range("J55").Value = 651.1539
$ws.Range("K55").Value = 167.33333
$ws.Range("L55").Value = 651.1539
$ws.Range("M55").Value = 5.666670000000011
$ws.Range("N55").Value = -997.1539
$ws.Range("H68").Value = 1140
$ws.Range("I68").Value = 1140
$ws.Range("K68").Value = 1140
$ws.Range("M68").Value = -391
$ws.Range("H71").Value = 1140
$ws.Range("I71").Value = 1140
$ws.Range("K71").Value = 5700
$ws.Range("M71").Value = -1956
$ws.Range("H99").Value = 23450.8
$ws.Range("I99").Value = 23450.8
$ws.Range("K99").Value = 23450.8
$ws.Range("M99").Value = -20455.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 33925
$ws.Range("I87").Value = 20295
$ws.Range("J87").Value = 47555
$ws.Range("K87").Value = 20295
$ws.Range("L87").Value = 47555
$ws.Range("M87").Value = -19047
$ws.Range("N87").Value = -50051
$ws.Range("H90").Value = 33925
$ws.Range("I90").Value = 20295
$ws.Range("J90").Value = 47555
$ws.Range("K90").Value = 60885
$ws.Range("L90").Value = 142665
$ws.Range("M90").Value = -54645
$ws.Range("N90").Value = -155145
$ws.Range("H122").Value = 2833.1667
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3499.5
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 10498.5
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -15398.5
